$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$vals = @(1.533923325814044,1.39310769697704,1.306397238246859,1.271001544175704,1.265120527843578,1.305920121524139,1.485423034714415,1.835375289934859,2.091159653628495,2.207221420493795,2.251126832654052,2.24167303611398,2.21083445221268,2.191939050922485,2.083568609095948,2.017009561521775,1.978698837251329,1.965722788667449,2.02409777024269,2.219893720573907,2.347596428553231,2.279463684251652,2.020893331211312,1.740931530175601)
$rng = $ws.Range("B2:B25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(0.006819927177946994,0.006082749589296554,0.00562887778522736,0.005443627510953064,0.005412849644660156,0.005626380600133984,0.006566016825079402,0.008398076470875537,0.009736839527960228,0.01034414423066465,0.01057385375043651,0.01052439370776881,0.01036304795813692,0.01026418416033437,0.009697114845344856,0.009348786246931695,0.009148277415796713,0.009080361751927057,0.009385883055415434,0.01041044641402067,0.01107851307956054,0.01072210065898105,0.009369112366194088,0.007903671744422525)
$rng = $ws.Range("D2:D25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(0.6256421266332808,0.5840387343000941,0.5585696683258448,0.5482114612186422,0.546492784095733,0.5584298881536256,0.6112825530356503,0.7154600416311467,0.7922440681288094,0.8272119715919501,0.8404574359122279,0.8376046418238872,0.8283016155927214,0.8226036997109674,0.7899594656907425,0.769941955289525,0.7584321308818147,0.7545357913617039,0.7720724806627146,0.8310340477366651,0.8695906909629372,0.8490108269818535,0.771109274659949,0.6872285346630065)
$rng = $ws.Range("E2:E25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(0.4829650877107099,0.4685925045170052,0.4602603782579351,0.4569881227318788,0.4564521837111997,0.4602157498134716,0.4779067115458275,0.5165436681299411,0.5473880330635836,0.5619651495510851,0.567564472394892,0.5663550192196851,0.5624242156905552,0.5600268337616683,0.5464464359813661,0.5382556730409789,0.5335958945642716,0.532026965786244,0.5391222753615068,0.5635766298944986,0.58002146534659,0.5712019671964157,0.5387303312562324,0.5056635425087137)
$rng = $ws.Range("F2:F25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(0.002364016048753207,0.002367478910217889,0.002369718658227406,0.002370660012201982,0.002370818055627005,0.002369731237365192,0.002365186532596751,0.002357171177671451,0.002351823306110737,0.002349506680731002,0.002348646044771669,0.002348830660186136,0.002349435542956088,0.002349808213238491,0.002351977033093142,0.002353337223540723,0.002354130506374896,0.00235440097945853,0.002353191297545448,0.002349257423775554,0.002346783244825905,0.002348094927520185,0.002353257235635849,0.002359244126534825)
$rng = $ws.Range("G2:G25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(0.7974680776926135,0.7922380670662861,0.7898031826238565,0.7890063475565299,0.7888858410223989,0.7897916448127731,0.795503730030859,0.8128610603524109,0.8293661198542281,0.8376902153685819,0.8409596631423284,0.8402503126071679,0.837956844146504,0.8365673031334069,0.8288385351186918,0.8243061225265222,0.8217759916371961,0.8209325232041351,0.8247806580563406,0.8386273082245026,0.8483605729774126,0.8431031931419284,0.8245658850070612,0.8075065089843818)
$rng = $ws.Range("I2:I25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(0.5386825466486869,0.4779951572535026,0.4406296734420891,0.4253779117987051,0.422843884501674,0.4404240829345838,0.5177794800370634,0.6686219202977952,0.7788948607842769,0.8289354134146265,0.8478660490869458,0.8437898413544644,0.8304932269277572,0.8223462125909293,0.7756220351914749,0.7469260575164185,0.7304093518518471,0.7248151258671953,0.7499819942561032,0.8343992779528264,0.8894618131125753,0.8600841802238506,0.7486004643934905,0.6279095568031039)
$rng = $ws.Range("L2:L25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @(1.511383579394447,1.477306341264352,1.458037223939186,1.450598296548208,1.449387959291357,1.457935229638025,1.499288772734332,1.593634836367357,1.671210754981075,1.708335001565246,1.72265978332149,1.719562781624319,1.709508150038147,1.703384205469149,1.668821756096975,1.648090702733157,1.636339259043012,1.632389964708523,1.650279684141566,1.71245418015198,1.754644271144116,1.731983367632722,1.649289525303203,1.566674280790323)
$rng = $ws.Range("O2:O25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}
